$d = $word.ActiveDocument
$vt = [char]11
$vt2 = @($vt, $vt) -join ""

# --- Short Summary sentence ---
$find1 = 'The summary of the article discusses the impact and potential of Artificial Intelligence AI in various sectors, including education, business, and media.'
$repl1 = 'The summary of the text discusses the growing impact and use of Artificial Intelligence AI in various sectors, including education, economy, media, and healthcare.'
$found1 = $d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $repl1, 2)
Write-Host "Change 1 (Short Summary sentence) found:" $found1

# --- Summary numbered list ---
$find2 = @('1. The global market value of AI in 2022 was estimated to be $119.78 billion.', '2. AI is expected to contribute $15.7 trillion to the global economy by 2030, representing a growth of 13x over the next 8 years.', '3. There are currently around 97 million people working in the AI industry globally and this number is projected to grow significantly.', '4. ChatGPT, a large-language model, is trained on billions of words of text data from various sources, enabling it to predict words and sentences in certain sequences based on patterns it has learned. However, it does not understand things in the conventional sense and its output should be taken with a pinch of salt.', '5. AI has the potential to revolutionize many industries, including entertainment, medical diagnosis, and legal services, but it also raises concerns about job displacement and ethical issues such as bias and privacy.', '6. Trust is key to the safe expansion of the use of AI solutions around the world. While there are tasks that can be better suited to automation with technology, the use of AI is still contingent upon human intelligence and awareness.', '7. The age of AI remains fraught with anxiety, with some people believing it will completely replace humans while others believe we are many generations away from this happening.', '8. In the UAE, Minister of Education Ahmed Belhoul Al-Falasi has emphasized the importance of using AI as a tool to assist in content creation and improve efficiency, rather than demonizing it as a threat to human creativity and employment.', '9. Jenna Burrell, director of research at Data Society, an independent non-profit research organization based in California, cautions that while AI can be useful for journalists, its information should be fact-checked as it is not up-to-date.', '10. Large-language models such as ChatGPT are limited to generating text based on patterns and do not understand things in the same way humans do. They are also not necessarily the most sophisticated form of AI, with other models such as reinforcement learning, generative adversarial networks, and symbolic AI emerging as alternatives.') -join $vt
$repl2 = @('1. The global market value of AI in 2022 is estimated to be $119.78 billion.', '2. AI is expected to contribute $15.7 trillion to the global economy by 2030, marking a projected growth of 13x over the next 8 years.', '3. The number of people working in the AI industry is projected to be 97 million by 2024.', 'The text also discusses the use of ChatGPT, an AI model developed by OpenAI, and its capabilities, limitations, and potential applications. Some examples include content creation on social media, blogs, and websites, writing business plans and reports, emails and presentations, legal documents such as contracts and medical summaries, and responding to customer inquiries and complaints.', 'However, the text also emphasizes that ChatGPT does not possess creativity, emotion, or personal perspective like human writers and is limited to generating text based on patterns it has learned during its training on text data. There are concerns about the use of AI in academic settings, with some institutions taking measures to prevent cheating using AI-detection programs.', 'Overall, while there is significant investment potential in the AI industry and many possible applications for AI technology, there is also anxiety and concern over its safe expansion and the role it will play in society. Trust is considered key to the safe and ethical use of AI solutions around the world.') -join $vt2
$found2 = $d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $repl2, 2)
Write-Host "Change 2 (Summary numbered list) found:" $found2

# --- Question 1 ---
$find3 = 'In the analyzed text, the media framing of the public discussion about ChatGPT can be observed through the use of several metaphors and recurring themes. 1. **Tool metaphor** ChatGPT is often described as a tool or instrument, such as I even passed exams at Minnesota University Law School and assist in content creation on social media, blogs, and websites. This framing emphasizes that ChatGPT is not an independent entity but rather a technology designed to perform specific tasks for human users. 2. **Pattern-matching metaphor** The text describes ChatGPT as a machine that matches patterns, which highlights its ability to generate text based on the data it was trained on without understanding in a conventional sense. This framing focuses on ChatGPT s limitations and helps set realistic expectations for users about its capabilities. 3. **Potential-threat metaphor** Although not overly prominent, there is an underlying suggestion that ChatGPT may pose a potential threat to certain professions, such as journalism, writing, or even higher education. This is evident in statements like There is a need for fact-checking and I am not able to replicate human innovation, creativity, skepticism, and reasoning.  4. **Caution metaphor** There are recurring cautions about the use of AI, with phrases such as need to take ChatGPT s answers with a pinch of salt or emphasizing the importance of fact-checking when using it for information gathering. This framing aims to educate users on the potential pitfalls and limitations of relying solely on AI. 5. **Promise metaphor** The text highlights the immense investment potential in AI, with forecasted values in the trillions of dollars, as well as its many possible applications across various industries like entertainment, medical diagnosis, etc. This framing focuses on the exciting possibilities that ChatGPT and other AI solutions could bring to society and the economy. Overall, these metaphors help shape public perception about ChatGPT by emphasizing its limitations, potential threats, and promises while encouraging caution when using it for various tasks.'
$parts3 = @('In the analyzed text, the media framing of the public discussion about ChatGPT can be observed through several recurring metaphors. Here are some key metaphors and their implications:', '1. **Pattern Matching**: "LLMs do not understand things in a conventional sense and they are only as good, or as accurate, as the information with which they are provided. They are essentially machines for matching patterns." This metaphor emphasizes that ChatGPT is not intelligent in a human-like manner but rather relies on patterns it has learned from data.', '2. **Creation and Innovation**: "Human writers bring creativity, emotion and personal perspective that I am not able to replicate" - this metaphor highlights the unique ability of humans to create original content and innovate compared to AI.', '3. **Dependency and Reliance**: "People need to take ChatGPT''s answers with a pinch of salt." This metaphor suggests that while AI can be useful, it should not be blindly trusted as its information may not always be accurate or up-to-date.', '4. **Replacement Fear**: "Don''t demonize AI as it will be a part of our lives... I think we are many generations away from when AI becomes greater than human capabilities." This metaphor conveys the anxiety surrounding AI''s potential to replace human jobs and abilities, while also acknowledging that humans still possess unique qualities unattainable by AI.', '5. **Promise and Potential**: "With forecasts valuing in the trillions of dollars" - this metaphor underscores the economic potential and investment opportunities associated with AI.', '6. **Age and Development**: "I think were many generations away from when AI becomes greater than human capabilities." This metaphor positions AI as a developing entity, growing and evolving over time.', '7. **Cautious Optimism**: "There''s as much optimism as there is pessimism over AI." This metaphor reflects a balanced perspective on the future of AI, acknowledging both its benefits and potential risks.')
$breaks3 = @($vt2, $vt, $vt, $vt, $vt, $vt, $vt)
$repl3 = $parts3[0]
for ($j = 0; $j -lt $breaks3.Length; $j++) {
    $repl3 = $repl3 + $breaks3[$j] + $parts3[$j + 1]
}
$found3 = $d.Content.Find.Execute($find3, $true, $false, $false, $false, $false, $true, 1, $false, $repl3, 2)
Write-Host "Change 3 (Question 1) found:" $found3

# --- Question 2 ---
$find4 = @('The provided text mainly covers the following perspectives and aspects of Artificial Intelligence AI  1. Application and use cases of AI across various sectors such as education, media, law, and business. 2. Investment potential in the AI industry with forecasted valuations in trillions of dollars. 3. The limitations of AI, such as its inability to understand things in a conventional sense or replicate human creativity and emotion. 4. Concerns about trust, ethics, and the responsible use of AI in various aspects of society. 5. The potential impact of AI on employment, particularly in relation to the number of people projected to work in the industry by 2024 97 million . 6. Efforts to prevent cheating using AI, such as programs that claim to catch text written by AI. 7. Discussions around fact-checking and the need for human oversight when using AI for news coverage or other critical information dissemination. 8. The role of governments and educational institutions in regulating and promoting responsible AI adoption. Aspects that seem to be somewhat ignored in this text include  1. A comprehensive exploration of alternative models of AI e.g., reinforcement learning, generative adversarial networks, symbolic AI and their current status or potential impact on the landscape. 2. The social and cultural implications of widespread AI adoption beyond employment and education, such as privacy concerns, bias in AI systems, or the psychological effects of increased reliance on machines for everyday tasks. 3. Detailed discussion of the progress and challenges in developing advanced forms of AI, such as AGI artificial general intelligence or ASI artificial superintelligence . 4. Exploration of potential collaborations between humans and AI in various sectors to create a symbiotic relationship rather than a competitive one.') -join $vt
$repl4 = @('The provided text discusses a wide range of perspectives and aspects related to Artificial Intelligence (AI), including its current market value, expected growth, potential economic contribution, and its various applications in industries such as education, media, law, and business. The text also touches upon the ethical considerations surrounding AI, such as preventing cheating and ensuring accuracy and reliability of information provided by AI.', 'One aspect that seems to be widely covered is the potential for AI to replace certain professional functions, particularly those involving content creation and analysis. However, there is also an emphasis on AI''s limitations, specifically its inability to fully imitate human creativity, emotion, skepticism, and reasoning. This suggests that while AI is seen as a powerful tool, it is not expected to completely replace humans across all industries.', 'Another aspect that is covered extensively is the potential economic impact of AI. The text mentions that AI is expected to contribute significantly to the global economy by 2030, with trillions of dollars forecasted in investment and growth. This suggests a focus on AI''s financial potential and its role in driving economic development.', 'However, there are some aspects that appear to be overlooked or under-discussed in this text. For instance, there is limited discussion about the social implications of AI, such as job displacement due to automation, privacy concerns, and potential AI-related ethical dilemmas. Additionally, while the text touches upon the need for fact-checking AI output, there is no substantial exploration of strategies for ensuring accuracy in AI-generated information or the role of humans in validating AI data. Lastly, there is no mention of the environmental impact of AI, which could be significant given the energy consumption required to power large-scale AI systems.', 'In summary, while the text provides a comprehensive overview of AI''s current status and potential future developments, it seems to underrepresent discussions on social implications, environmental impact, and strategies for validating AI output.') -join $vt2
$found4 = $d.Content.Find.Execute($find4, $true, $false, $false, $false, $false, $true, 1, $false, $repl4, 2)
Write-Host "Change 4 (Question 2) found:" $found4

# --- Question 3 ---
$find5 = 'Not mentioned'
$repl5 = 'Not mentioned. The article primarily focuses on global discussions and developments regarding AI, and there is no specific mention of how the Arabic world leverages AI in the context provided.'
$found5 = $d.Content.Find.Execute($find5, $true, $false, $false, $false, $false, $true, 1, $false, $repl5, 2)
Write-Host "Change 5 (Question 3) found:" $found5

# --- Question 4 ---
$find6 = @('The article highlights the growing influence and potential of AI, with investments in trillions and applications ranging from entertainment to medical diagnosis. However, concerns about its impact on employment, ethics, and trust remain. The final message emphasizes the importance of human intelligence and awareness in the safe expansion of AI use.') -join $vt
$repl6 = @('1) The rapid advancement and integration of AI into various sectors present immense opportunities for investment, economic growth, and productivity improvements.', '2) However, the use of AI also raises concerns about job displacement, ethics, privacy, and the potential for misinformation or bias.', '3) It is crucial to foster trust in AI solutions by ensuring transparency, accountability, and ethical guidelines for their development and deployment.', '4) Human intelligence and awareness remain vital in the safe expansion of AI usage around the world.') -join $vt
$found6 = $d.Content.Find.Execute($find6, $true, $false, $false, $false, $false, $true, 1, $false, $repl6, 2)
Write-Host "Change 6 (Question 4) found:" $found6

# --- Entities list ---
$find7 = 'AFP However, Generative Pre-Trained Transformer, Google, DeenSquare, LEAP, Elon Musk, YouTube, Bard, PwC Middle Easts, AI, Guardian, Sam Altman, Intelligence, Googles, Scott Nowson, UAE, Alex Hern, Reid Ho, Noaman Sayed, Nowson, Dubais Museum of the Future, Edge, Ahmed Belhoul Al-Falasi, OpenAI, Burrell, Arab News, the World Economic Forum, AFP Founded, Microsoft, Alexa, Omar Sultan Al-Olama, Bing, the Sciences Po school, Minnesota University Law School, Dan Milmo, ChatGPT, Peter Thiel, James Webb, New York Citys, GPT, National Strategy for Data, Jenna Burrell, LinkedIn, Data Society'
$repl7 = 'Minnesota University Law School, OpenAI, Noaman Sayed, Dan Milmo, AFP Founded, Nowson, Bard, UAE, Data Society, New York Citys, Alex Hern, Edge, the Sciences Po school, AI, Sam Altman, Intelligence, Generative Pre-Trained Transformer, Elon Musk, LinkedIn, Burrell, PwC Middle Easts, Googles, Bing, GPT, Peter Thiel, Dubais Museum of the Future, the World Economic Forum, National Strategy for Data, LEAP, AFP However, Omar Sultan Al-Olama, Guardian, ChatGPT, YouTube, Alexa, Reid Ho, DeenSquare, Ahmed Belhoul Al-Falasi, Microsoft, Scott Nowson, James Webb, Jenna Burrell, Arab News, Google'
$found7 = $d.Content.Find.Execute($find7, $true, $false, $false, $false, $false, $true, 1, $false, $repl7, 2)
Write-Host "Change 7 (Entities list) found:" $found7

